$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D so numeric-looking text (e.g. '545.13') is written as
# text, not auto-converted to a number -- matches source data which stores
# prices as inline strings.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = '63.494.27'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '3.080.48'
$ws.Range("E3").Value = '  -0.44%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '545.13'
$ws.Range("E5").Value = '  -0.85%  '

$ws.Range("D6").Value = '139.58'
$ws.Range("E6").Value = '  +1.62%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '3.075.90'
$ws.Range("E8").Value = '  -0.30%  '

$ws.Range("E9").Value = '  +0.45%  '

$ws.Range("E10").Value = '  +0.54%  '

$ws.Range("D11").Value = '6.42'
$ws.Range("E11").Value = '  +2.64%  '

$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  -2.78%  '

$ws.Range("E13").Value = '  +3.58%  '

$ws.Range("D14").Value = '35.01'
$ws.Range("E14").Value = '  -1.57%  '

$ws.Range("D15").Value = '3.580.47'
$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("D16").Value = '63.502.60'
$ws.Range("E16").Value = '  +0.27%  '

$ws.Range("E17").Value = '  +0.94%  '

$ws.Range("D18").Value = '3.081.37'
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("D19").Value = '6.65'
$ws.Range("E19").Value = '  -1.16%  '

$ws.Range("D20").Value = '475.96'
$ws.Range("E20").Value = '  -2.74%  '

$ws.Range("D21").Value = '13.48'
$ws.Range("E21").Value = '  -1.36%  '

$ws.Range("D22").Value = '0.700'
$ws.Range("E22").Value = '  -2.62%  '

$ws.Range("D23").Value = '7.09'
$ws.Range("E23").Value = '  -2.26%  '

$ws.Range("D24").Value = '78.69'
$ws.Range("E24").Value = '  -0.54%  '

$ws.Range("D25").Value = '12.24'
$ws.Range("E25").Value = '  -1.25%  '

$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("E27").Value = '  -1.24%  '

$ws.Range("D28").Value = '8.01'
$ws.Range("E28").Value = '  -5.68%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").Value = '26.23'
$ws.Range("E30").Value = '  -1.43%  '

$ws.Range("E31").Value = '  -4.04%  '

$ws.Range("E32").Value = '  +1.64%  '

$ws.Range("D33").Value = '58.77'
$ws.Range("E33").Value = '  +1.12%  '

$ws.Range("E34").Value = '  -7.54%  '

$ws.Range("D35").Value = '5.50'
$ws.Range("E35").Value = '  +6.97%  '

$ws.Range("D36").Value = '490.86'
$ws.Range("E36").Value = '  -4.31%  '

$ws.Range("D37").Value = '6.01'
$ws.Range("E37").Value = '  -0.31%  '

$ws.Range("D38").Value = '3.265.37'
$ws.Range("E38").Value = '  +3.43%  '

$ws.Range("E39").Value = '  +0.46%  '

$ws.Range("D40").Value = '0.0799'
$ws.Range("E40").Value = '  -0.61%  '

$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("D42").Value = '8.15'
$ws.Range("E42").Value = '  -0.29%  '

$ws.Range("D43").Value = '2.60'
$ws.Range("E43").Value = '  -2.35%  '

$ws.Range("D44").Value = '0.254'
$ws.Range("E44").Value = '  -2.41%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").Value = '25.71'
$ws.Range("E46").Value = '  +1.95%  '

$ws.Range("D47").Value = '123.69'
$ws.Range("E47").Value = '  +2.21%  '

$ws.Range("E48").Value = '  -1.74%  '

$ws.Range("D49").Value = '0.0₃0530'
$ws.Range("E49").Value = '  +5.16%  '

$ws.Range("E50").Value = '  +0.57%  '

$ws.Range("D51").Value = '2.02'
$ws.Range("E51").Value = '  -0.83%  '

# Restore default cell style now that the text values are committed.
$priceCol.Style = "Normal"
